$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.5444614721754988
$ws.Cells.Item(2, 3).Value = 0.2199998318697993
$ws.Cells.Item(2, 4).Value = 0.07825835179046692
$ws.Cells.Item(2, 5).Value = 0.1354707382770037
$ws.Cells.Item(2, 7).Value = 0.7963246162988327
$ws.Cells.Item(2, 8).Value = 0.8801369325400543
$ws.Cells.Item(2, 11).Value = 0.2744477690766303
$ws.Cells.Item(2, 12).Value = 0.1902934636307378
$ws.Cells.Item(2, 13).Value = 0.1568398313437918
$ws.Cells.Item(2, 14).Value = 1.929792603485881
$ws.Cells.Item(2, 15).Value = 3.369760866746361
$ws.Cells.Item(3, 2).Value = 0.5088688634830305
$ws.Cells.Item(3, 3).Value = 0.2197934695591854
$ws.Cells.Item(3, 4).Value = 0.07103394885483283
$ws.Cells.Item(3, 5).Value = 0.1359704241426751
$ws.Cells.Item(3, 7).Value = 0.7980827652367921
$ws.Cells.Item(3, 8).Value = 0.8846671997178461
$ws.Cells.Item(3, 11).Value = 0.2422383370241619
$ws.Cells.Item(3, 12).Value = 0.1876728200020494
$ws.Cells.Item(3, 13).Value = 0.1501696946828055
$ws.Cells.Item(3, 14).Value = 1.948076651629993
$ws.Cells.Item(3, 15).Value = 3.382704102399714
$ws.Cells.Item(4, 2).Value = 0.4871974121235212
$ws.Cells.Item(4, 3).Value = 0.2196798228123953
$ws.Cells.Item(4, 4).Value = 0.0666308899424024
$ws.Cells.Item(4, 5).Value = 0.1363276894668459
$ws.Cells.Item(4, 7).Value = 0.7995943919714463
$ws.Cells.Item(4, 8).Value = 0.8877772601946958
$ws.Cells.Item(4, 11).Value = 0.2224574114253528
$ws.Cells.Item(4, 12).Value = 0.1861536171497136
$ws.Cells.Item(4, 13).Value = 0.1461408468812451
$ws.Cells.Item(4, 14).Value = 1.959877133482028
$ws.Cells.Item(4, 15).Value = 3.392245412146607
$ws.Cells.Item(5, 2).Value = 0.4784126276987877
$ws.Cells.Item(5, 3).Value = 0.2196368180707431
$ws.Cells.Item(5, 4).Value = 0.06484487948836204
$ws.Cells.Item(5, 5).Value = 0.1364859881171157
$ws.Cells.Item(5, 7).Value = 0.8003190875331541
$ws.Cells.Item(5, 8).Value = 0.8891273291805106
$ws.Cells.Item(5, 11).Value = 0.2143959096138985
$ws.Cells.Item(5, 12).Value = 0.1855572014470255
$ws.Cells.Item(5, 13).Value = 0.1445159388752764
$ws.Cells.Item(5, 14).Value = 1.964830296744615
$ws.Cells.Item(5, 15).Value = 3.396534706336254
$ws.Cells.Item(6, 2).Value = 0.4769567476513146
$ws.Cells.Item(6, 3).Value = 0.2196298775952386
$ws.Cells.Item(6, 4).Value = 0.06454881445300487
$ws.Cells.Item(6, 5).Value = 0.1365130418098364
$ws.Cells.Item(6, 7).Value = 0.8004459886859294
$ws.Cells.Item(6, 8).Value = 0.8893565045911629
$ws.Cells.Item(6, 11).Value = 0.2130572793357572
$ws.Cells.Item(6, 12).Value = 0.1854595384529603
$ws.Cells.Item(6, 13).Value = 0.1442471472561628
$ws.Cells.Item(6, 14).Value = 1.965661486411041
$ws.Cells.Item(6, 15).Value = 3.397271176098883
$ws.Cells.Item(7, 2).Value = 0.4870787483405934
$ws.Cells.Item(7, 3).Value = 0.2196792294127334
$ws.Cells.Item(7, 4).Value = 0.06660676964121137
$ws.Cells.Item(7, 5).Value = 0.136329772845011
$ws.Cells.Item(7, 7).Value = 0.7996037253177875
$ws.Cells.Item(7, 8).Value = 0.8877951327383755
$ws.Cells.Item(7, 11).Value = 0.2223486930916607
$ws.Cells.Item(7, 12).Value = 0.1861454817883086
$ws.Cells.Item(7, 13).Value = 0.1461188642942837
$ws.Cells.Item(7, 14).Value = 1.959943348976077
$ws.Cells.Item(7, 15).Value = 3.392301634524287
$ws.Cells.Item(8, 2).Value = 0.5321516017701242
$ws.Cells.Item(8, 3).Value = 0.2199259810867886
$ws.Cells.Item(8, 4).Value = 0.07576059887624353
$ws.Cells.Item(8, 5).Value = 0.1356325711650452
$ws.Cells.Item(8, 7).Value = 0.7968411604516348
$ws.Cells.Item(8, 8).Value = 0.8816308570558817
$ws.Cells.Item(8, 11).Value = 0.2633430856475911
$ws.Cells.Item(8, 12).Value = 0.1893712512706145
$ws.Cells.Item(8, 13).Value = 0.1545262074526015
$ws.Cells.Item(8, 14).Value = 1.935977830978549
$ws.Cells.Item(8, 15).Value = 3.373893000694778
$ws.Cells.Item(9, 2).Value = 0.6219659941642988
$ws.Cells.Item(9, 3).Value = 0.220512533938404
$ws.Cells.Item(9, 4).Value = 0.09397048236958483
$ws.Cells.Item(9, 5).Value = 0.1346647750037953
$ws.Cells.Item(9, 7).Value = 0.7948516560993824
$ws.Cells.Item(9, 8).Value = 0.8721447133119824
$ws.Cells.Item(9, 11).Value = 0.3436839891006684
$ws.Cells.Item(9, 12).Value = 0.1964076002063777
$ws.Cells.Item(9, 13).Value = 0.1715373444769277
$ws.Cells.Item(9, 14).Value = 1.893533753772555
$ws.Cells.Item(9, 15).Value = 3.350432183885744
$ws.Cells.Item(10, 2).Value = 0.6887999268325871
$ws.Cells.Item(10, 3).Value = 0.2210048650057246
$ws.Cells.Item(10, 4).Value = 0.1075078553106152
$ws.Cells.Item(10, 5).Value = 0.1341960572877277
$ws.Cells.Item(10, 7).Value = 0.7954798751518268
$ws.Cells.Item(10, 8).Value = 0.8667562024485136
$ws.Cells.Item(10, 11).Value = 0.4026652791149559
$ws.Cells.Item(10, 12).Value = 0.2020077526613875
$ws.Cells.Item(10, 13).Value = 0.1843505524163618
$ws.Cells.Item(10, 14).Value = 1.86512116343078
$ws.Cells.Item(10, 15).Value = 3.340890019656854
$ws.Cells.Item(11, 2).Value = 0.7193840373687408
$ws.Cells.Item(11, 3).Value = 0.2212419214583932
$ws.Cells.Item(11, 4).Value = 0.113700970676831
$ws.Cells.Item(11, 5).Value = 0.134035201398973
$ws.Cells.Item(11, 7).Value = 0.7962195941144898
$ws.Cells.Item(11, 8).Value = 0.86464703670174
$ws.Cells.Item(11, 11).Value = 0.4294846632501219
$ws.Cells.Item(11, 12).Value = 0.2046483753004793
$ws.Cells.Item(11, 13).Value = 0.1902471572144222
$ws.Cells.Item(11, 14).Value = 1.852796378929714
$ws.Cells.Item(11, 15).Value = 3.33821791053029
$ws.Cells.Item(12, 2).Value = 0.7309909493312432
$ws.Cells.Item(12, 3).Value = 0.2213335492330089
$ws.Cells.Item(12, 4).Value = 0.1160511408808134
$ws.Cells.Item(12, 5).Value = 0.1339817984863814
$ws.Cells.Item(12, 7).Value = 0.7965649702461661
$ws.Cells.Item(12, 8).Value = 0.8638974520476097
$ws.Cells.Item(12, 11).Value = 0.4396384327006331
$ws.Cells.Item(12, 12).Value = 0.2056616354539216
$ws.Cells.Item(12, 13).Value = 0.1924896932610096
$ws.Cells.Item(12, 14).Value = 1.848215589612973
$ws.Cells.Item(12, 15).Value = 3.337445795022887
$ws.Cells.Item(13, 2).Value = 0.7284900751683097
$ws.Cells.Item(13, 3).Value = 0.221313733211197
$ws.Cells.Item(13, 4).Value = 0.1155447693007972
$ws.Cells.Item(13, 5).Value = 0.1339929660879982
$ws.Cells.Item(13, 7).Value = 0.7964876850826812
$ws.Cells.Item(13, 8).Value = 0.8640567056105368
$ws.Cells.Item(13, 11).Value = 0.437451740261821
$ws.Cells.Item(13, 12).Value = 0.2054428207680701
$ws.Cells.Item(13, 13).Value = 0.1920062972540535
$ws.Cells.Item(13, 14).Value = 1.849198304181971
$ws.Cells.Item(13, 15).Value = 3.337601423942971
$ws.Cells.Item(14, 2).Value = 0.7203384402556878
$ws.Cells.Item(14, 3).Value = 0.2212494225933384
$ws.Cells.Item(14, 4).Value = 0.113894221219752
$ws.Cells.Item(14, 5).Value = 0.1340306575203343
$ws.Cells.Item(14, 7).Value = 0.79624670058638
$ws.Cells.Item(14, 8).Value = 0.8645843841402723
$ws.Cells.Item(14, 11).Value = 0.4303200658887931
$ws.Cells.Item(14, 12).Value = 0.204731470409925
$ws.Cells.Item(14, 13).Value = 0.1904314599642021
$ws.Cells.Item(14, 14).Value = 1.852417784122765
$ws.Cells.Item(14, 15).Value = 3.338149584228916
$ws.Cells.Item(15, 2).Value = 0.7153486132965554
$ws.Cells.Item(15, 3).Value = 0.2212102719808229
$ws.Cells.Item(15, 4).Value = 0.1128838589498287
$ws.Cells.Item(15, 5).Value = 0.1340547219911699
$ws.Cells.Item(15, 7).Value = 0.7961075890193996
$ws.Cells.Item(15, 8).Value = 0.8649139956983305
$ws.Cells.Item(15, 11).Value = 0.425951414795918
$ws.Cells.Item(15, 12).Value = 0.204297479471478
$ws.Cells.Item(15, 13).Value = 0.1894680755230453
$ws.Cells.Item(15, 14).Value = 1.85440105321503
$ws.Cells.Item(15, 15).Value = 3.338516565412874
$ws.Cells.Item(16, 2).Value = 0.6868047696877113
$ws.Cells.Item(16, 3).Value = 0.220989634289765
$ws.Cells.Item(16, 4).Value = 0.1071038189327425
$ws.Cells.Item(16, 5).Value = 0.1342076215891819
$ws.Cells.Item(16, 7).Value = 0.7954406658998892
$ws.Cells.Item(16, 8).Value = 0.8669009181797662
$ws.Cells.Item(16, 11).Value = 0.4009122955284568
$ws.Cells.Item(16, 12).Value = 0.2018370487353565
$ws.Cells.Item(16, 13).Value = 0.1839665487428661
$ws.Cells.Item(16, 14).Value = 1.86593869759607
$ws.Cells.Item(16, 15).Value = 3.341098209217279
$ws.Cells.Item(17, 2).Value = 0.6693399421750428
$ws.Cells.Item(17, 3).Value = 0.2208576170930172
$ws.Cells.Item(17, 4).Value = 0.1035668568360961
$ws.Cells.Item(17, 5).Value = 0.134314818199222
$ws.Cells.Item(17, 7).Value = 0.795147784158118
$ws.Cells.Item(17, 8).Value = 0.8682073908703103
$ws.Cells.Item(17, 11).Value = 0.3855483089936627
$ws.Cells.Item(17, 12).Value = 0.2003514444549239
$ws.Cells.Item(17, 13).Value = 0.1806088150914462
$ws.Cells.Item(17, 14).Value = 1.873170422635589
$ws.Cells.Item(17, 15).Value = 3.343109209303009
$ws.Cells.Item(18, 2).Value = 0.6593117130632606
$ws.Cells.Item(18, 3).Value = 0.2207829179970062
$ws.Cells.Item(18, 4).Value = 0.1015357772411818
$ws.Cells.Item(18, 5).Value = 0.134381405027197
$ws.Cells.Item(18, 7).Value = 0.795022054157414
$ws.Cells.Item(18, 8).Value = 0.8689910438820476
$ws.Cells.Item(18, 11).Value = 0.3767102974831289
$ws.Cells.Item(18, 12).Value = 0.1995057300958791
$ws.Cells.Item(18, 13).Value = 0.1786839261923774
$ws.Cells.Item(18, 14).Value = 1.877386415910928
$ws.Cells.Item(18, 15).Value = 3.34442295829075
$ws.Cells.Item(19, 2).Value = 0.6559192801278755
$ws.Cells.Item(19, 3).Value = 0.2207578387622533
$ws.Cells.Item(19, 4).Value = 0.1008486546056844
$ws.Cells.Item(19, 5).Value = 0.1344047975688625
$ws.Cells.Item(19, 7).Value = 0.7949868231420254
$ws.Cells.Item(19, 8).Value = 0.8692619089876672
$ws.Cells.Item(19, 11).Value = 0.3737177329073802
$ws.Cells.Item(19, 12).Value = 0.1992208934886435
$ws.Cells.Item(19, 13).Value = 0.1780332934958437
$ws.Cells.Item(19, 14).Value = 1.878823581383099
$ws.Cells.Item(19, 15).Value = 3.344894754936178
$ws.Cells.Item(20, 2).Value = 0.6711973401545777
$ws.Cells.Item(20, 3).Value = 0.2208715430319828
$ws.Cells.Item(20, 4).Value = 0.1039430325822224
$ws.Cells.Item(20, 5).Value = 0.1343028968311781
$ws.Cells.Item(20, 7).Value = 0.7951745398576548
$ws.Cells.Item(20, 8).Value = 0.8680649821757953
$ws.Cells.Item(20, 11).Value = 0.3871839441381724
$ws.Cells.Item(20, 12).Value = 0.2005086828461913
$ws.Cells.Item(20, 13).Value = 0.1809655911106347
$ws.Cells.Item(20, 14).Value = 1.872394745832252
$ws.Cells.Item(20, 15).Value = 3.342878879678409
$ws.Cells.Item(21, 2).Value = 0.7227320896528227
$ws.Cells.Item(21, 3).Value = 0.2212682618925896
$ws.Cells.Item(21, 4).Value = 0.1143788926409286
$ws.Cells.Item(21, 5).Value = 0.134019382984377
$ws.Cells.Item(21, 7).Value = 0.7963157125446401
$ws.Cells.Item(21, 8).Value = 0.864428059969498
$ws.Cells.Item(21, 11).Value = 0.4324148738096483
$ws.Cells.Item(21, 12).Value = 0.2049400504979815
$ws.Cells.Item(21, 13).Value = 0.1908937678223097
$ws.Cells.Item(21, 14).Value = 1.851469801397039
$ws.Cells.Item(21, 15).Value = 3.337982070988147
$ws.Cells.Item(22, 2).Value = 0.7565606155254727
$ws.Cells.Item(22, 3).Value = 0.221538369400875
$ws.Cells.Item(22, 4).Value = 0.1212282878852733
$ws.Cells.Item(22, 5).Value = 0.1338778532129119
$ws.Cells.Item(22, 7).Value = 0.7974419161679691
$ws.Cells.Item(22, 8).Value = 0.8623373375800014
$ws.Cells.Item(22, 11).Value = 0.4619631430146569
$ws.Cells.Item(22, 12).Value = 0.2079137741734485
$ws.Cells.Item(22, 13).Value = 0.1974384154594446
$ws.Cells.Item(22, 14).Value = 1.838297403991625
$ws.Cells.Item(22, 15).Value = 3.336179102118138
$ws.Cells.Item(23, 2).Value = 0.7384924077245785
$ws.Cells.Item(23, 3).Value = 0.2213932246886898
$ws.Cells.Item(23, 4).Value = 0.1175700026710871
$ws.Cells.Item(23, 5).Value = 0.1339493927638244
$ws.Cells.Item(23, 7).Value = 0.7968060403795363
$ws.Cells.Item(23, 8).Value = 0.8634270336945349
$ws.Cells.Item(23, 11).Value = 0.4461940069213028
$ws.Cells.Item(23, 12).Value = 0.2063195681384968
$ws.Cells.Item(23, 13).Value = 0.1939403324338613
$ws.Cells.Item(23, 14).Value = 1.845281703015926
$ws.Cells.Item(23, 15).Value = 3.337013585975825
$ws.Cells.Item(24, 2).Value = 0.6703575714404053
$ws.Cells.Item(24, 3).Value = 0.2208652433775242
$ws.Cells.Item(24, 4).Value = 0.1037729561520848
$ws.Cells.Item(24, 5).Value = 0.134308271036474
$ws.Cells.Item(24, 7).Value = 0.7951623107479833
$ws.Cells.Item(24, 8).Value = 0.8681292637659794
$ws.Cells.Item(24, 11).Value = 0.3864444891854646
$ws.Cells.Item(24, 12).Value = 0.2004375692636415
$ws.Cells.Item(24, 13).Value = 0.1808042754781241
$ws.Cells.Item(24, 14).Value = 1.872745247470395
$ws.Cells.Item(24, 15).Value = 3.342982520777497
$ws.Cells.Item(25, 2).Value = 0.5975184035944494
$ws.Cells.Item(25, 3).Value = 0.2203429886685626
$ws.Cells.Item(25, 4).Value = 0.08901641402144378
$ws.Cells.Item(25, 5).Value = 0.134883950675535
$ws.Cells.Item(25, 7).Value = 0.7950228716691328
$ws.Cells.Item(25, 8).Value = 0.8744329539309348
$ws.Cells.Item(25, 11).Value = 0.3219564700871445
$ws.Cells.Item(25, 12).Value = 0.1944282793781795
$ws.Cells.Item(25, 13).Value = 0.1668797130906405
$ws.Cells.Item(25, 14).Value = 1.904529086234155
$ws.Cells.Item(25, 15).Value = 3.355426982095111
